$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 635.5833
$ws.Range("J17").Value = 635.5833
$ws.Range("L17").Value = 1906.7499
$ws.Range("N17").Value = -2242.7499

$ws.Range("H18").Value = 1180.5454
$ws.Range("I18").Value = 948.6
$ws.Range("J18").Value = 3500
$ws.Range("K18").Value = 948.6
$ws.Range("L18").Value = 3500
$ws.Range("M18").Value = -664.6
$ws.Range("N18").Value = -4068

$ws.Range("H98").Value = 1289.1333
$ws.Range("I98").Value = 997.5454999999999
$ws.Range("K98").Value = 997.5454999999999
$ws.Range("M98").Value = 500.4545000000001

$ws.Range("H106").Value = 2766.2727
$ws.Range("I106").Value = 2766.2727
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2766.2727
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -2135.2727

$ws.Range("H122").Value = 1289.1333
$ws.Range("I122").Value = 997.5454999999999
$ws.Range("K122").Value = 2992.6365
$ws.Range("M122").Value = -542.6364999999996

$ws.Range("H138").Value = 1644.2808
$ws.Range("I138").Value = 1318.8889
$ws.Range("J138").Value = 1794.4615
$ws.Range("K138").Value = 3956.6667
$ws.Range("L138").Value = 5383.3845
$ws.Range("M138").Value = 1183.3333
$ws.Range("N138").Value = -15663.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 843.25
$ws.Range("J2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("N2").Value = -1226

$ws.Range("H45").Value = 11366331
$ws.Range("I45").Value = 3310.1667
$ws.Range("K45").Value = 3310.1667
$ws.Range("M45").Value = -2933.1667

$ws.Range("H97").Value = 936.0909
$ws.Range("I97").Value = 729.7
$ws.Range("K97").Value = 729.7
$ws.Range("M97").Value = -233.7

$ws.Range("H116").Value = 843.25
$ws.Range("J116").Value = 1000
$ws.Range("L116").Value = 1000
$ws.Range("N116").Value = -5588

$ws.Range("H122").Value = 1749.6207
$ws.Range("I122").Value = 1519.7693
$ws.Range("K122").Value = 4559.3079
$ws.Range("M122").Value = -2109.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 843.25
$ws.Range("J3").Value = 1000
$ws.Range("L3").Value = 1000
$ws.Range("N3").Value = -1228

$ws.Range("H62").Value = 60000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 60000
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = 60000
$ws.Range("N62").Value = -61372

$ws.Range("H65").Value = 60000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 60000
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = 180000
$ws.Range("N65").Value = -186864

$ws.Range("H86").Value = 5376.5386
$ws.Range("I86").Value = 3110.375
$ws.Range("J86").Value = 9002.4
$ws.Range("K86").Value = 3110.375
$ws.Range("L86").Value = 9002.4
$ws.Range("M86").Value = -1987.375
$ws.Range("N86").Value = -11248.4

$ws.Range("H89").Value = 5376.5386
$ws.Range("I89").Value = 3110.375
$ws.Range("J89").Value = 9002.4
$ws.Range("K89").Value = 15551.875
$ws.Range("L89").Value = 45012
$ws.Range("M89").Value = -9935.875
$ws.Range("N89").Value = -56244

$ws.Range("H105").Value = 74446.07000000001
$ws.Range("I105").Value = 145302.42
$ws.Range("K105").Value = 145302.42
$ws.Range("M105").Value = -143555.42

$ws.Range("H107").Value = 1342.4706
$ws.Range("I107").Value = 1058.4286
$ws.Range("J107").Value = 2668
$ws.Range("K107").Value = 1058.4286
$ws.Range("L107").Value = 2668
$ws.Range("M107").Value = 861.5714
$ws.Range("N107").Value = -6508

$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 614
$ws.Range("I22").Value = 684.5
$ws.Range("J22").Value = 379
$ws.Range("K22").Value = 684.5
$ws.Range("L22").Value = 379
$ws.Range("M22").Value = -334.5
$ws.Range("N22").Value = -1079

$ws.Range("H134").Value = 1870.8334
$ws.Range("I134").Value = 1400.125
$ws.Range("K134").Value = 4200.375
$ws.Range("M134").Value = -1665.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 6125.5
$ws.Range("I98").Value = 1003
$ws.Range("J98").Value = 7833
$ws.Range("K98").Value = 3009
$ws.Range("L98").Value = 23499
$ws.Range("M98").Value = -1511
$ws.Range("N98").Value = -26495

$ws.Range("H103").Value = 928.8333
$ws.Range("I103").Value = 914.5
$ws.Range("K103").Value = 2743.5
$ws.Range("M103").Value = -1864.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("N39").Value = 0

$ws.Range("H122").Value = 280228.47
$ws.Range("I122").Value = 359046.97
$ws.Range("J122").Value = 4363.75
$ws.Range("K122").Value = 1077140.91
$ws.Range("L122").Value = 13091.25
$ws.Range("M122").Value = -1074690.91
$ws.Range("N122").Value = -17991.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2142.8572
$ws.Range("J16").Value = 2333.3333
$ws.Range("L16").Value = 2333.3333
$ws.Range("N16").Value = -2673.3333

$ws.Range("H43").Value = 34000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H55").Value = 6096.9473
$ws.Range("I55").Value = 541.6429000000001
$ws.Range("J55").Value = 21651.8
$ws.Range("K55").Value = 541.6429000000001
$ws.Range("L55").Value = 21651.8
$ws.Range("M55").Value = -368.6429000000001
$ws.Range("N55").Value = -21997.8

$ws.Range("H61").Value = 2664.1667
$ws.Range("I61").Value = 999
$ws.Range("J61").Value = 5994.5
$ws.Range("K61").Value = 999
$ws.Range("L61").Value = 5994.5
$ws.Range("M61").Value = -797
$ws.Range("N61").Value = -6398.5

$ws.Range("H100").Value = 8645.137000000001
$ws.Range("I100").Value = 10328.23
$ws.Range("K100").Value = 10328.23
$ws.Range("M100").Value = -9787.23

$ws.Range("H113").Value = 2664.1667
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 5994.5
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 5994.5
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -10334.5

$ws.Range("H122").Value = 80004400
$ws.Range("I122").Value = 125004210
$ws.Range("J122").Value = 28576052
$ws.Range("K122").Value = 375012630
$ws.Range("L122").Value = 85728156
$ws.Range("M122").Value = -375010180
$ws.Range("N122").Value = -85733056

$ws.Range("H136").Value = 20835794
$ws.Range("J136").Value = 2262.0454
$ws.Range("L136").Value = 6786.1362
$ws.Range("N136").Value = -11886.1362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1622.1875
$ws.Range("I107").Value = 318.25
$ws.Range("J107").Value = 2926.125
$ws.Range("K107").Value = 954.75
$ws.Range("L107").Value = 8778.375
$ws.Range("M107").Value = 965.25
$ws.Range("N107").Value = -12618.375

$ws.Range("H122").Value = 2402.8462
$ws.Range("I122").Value = 2037.4445
$ws.Range("K122").Value = 6112.333500000001
$ws.Range("M122").Value = -3662.333500000001

$ws.Range("H136").Value = 13514841
$ws.Range("I136").Value = 17242526
$ws.Range("J136").Value = 1979
$ws.Range("L136").Value = 5937
$ws.Range("N136").Value = -11037
